$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New LexTALE participant rows appended below the existing data block
# (A1:B44 -> A1:B74). Every cell in this sheet is stored as text, even
# when the content is numeric-looking (e.g. "2400", "93.75"), and blank
# cells are still materialised (empty text) rather than left absent.
# Prefixing the literal with a leading apostrophe forces Excel to store
# the value as text instead of inferring a number; ClearFormats() then
# drops the quote-prefix formatting flag so no stray cell style is left
# behind, matching the plain (unstyled) cells used throughout the sheet.
$newRows = @(
    @{Row=45; A=""; B=""},
    @{Row=46; A="2400"; B="93.75"},
    @{Row=47; A="2500"; B="95"},
    @{Row=48; A="1008"; B="56.25"},
    @{Row=49; A="1009"; B="76.25"},
    @{Row=50; A="3001"; B="55"},
    @{Row=51; A="3003"; B="62.5"},
    @{Row=52; A="3005"; B="55"},
    @{Row=53; A="2200"; B="78,75"},
    @{Row=54; A="3004"; B="58.75"},
    @{Row=55; A=""; B=""},
    @{Row=56; A="3010"; B="48.5"},
    @{Row=57; A="3007"; B="78.75"},
    @{Row=58; A="2600"; B="57,5"},
    @{Row=59; A="3006"; B="80"},
    @{Row=60; A="3008"; B="55"},
    @{Row=61; A="3009"; B="73.75"},
    @{Row=62; A="2700"; B="91.25"},
    @{Row=63; A="4002"; B="78.75"},
    @{Row=64; A=""; B=""},
    @{Row=65; A="3002"; B="56.25"},
    @{Row=66; A="4003"; B="61.25"},
    @{Row=67; A="4004"; B="61.25"},
    @{Row=68; A="4005"; B="62.50"},
    @{Row=69; A="4006"; B="62.5"},
    @{Row=70; A="4008"; B="43.75"},
    @{Row=71; A="4007"; B="65"},
    @{Row=72; A="4007"; B="65"},
    @{Row=73; A="4007"; B="65"},
    @{Row=74; A="4009"; B="43.75"}
)

foreach ($item in $newRows) {
    $rowNum = $item.Row

    $aCell = $ws.Range("A" + $rowNum)
    $aCell.Value = "'" + $item.A
    $aCell.ClearFormats()

    $bCell = $ws.Range("B" + $rowNum)
    $bCell.Value = "'" + $item.B
    $bCell.ClearFormats()
}
